$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("new table with CFI etc")

# --- Row 2 (model 1): was all "NA" -> now has real numbers + "<0.001" ---
$ws.Range("C2").Value = 0.942
$ws.Range("D2").Value = 0.101
$ws.Range("E2").Value = 5129.226
$ws.Range("F2").Value = "<0.001"

# --- Row 3 (model 2): tweak CFI + AICc ---
$ws.Range("C3").Value = 0.929
$ws.Range("E3").Value = 5941.559

# --- Row 4 (model 3): now becomes "NA" across the board ---
$ws.Range("C4").Value = "NA"
$ws.Range("D4").Value = "NA"
$ws.Range("E4").Value = "NA"
$ws.Range("F4").Value = "NA"

# --- Row 6 (model 5): tweak numbers and bold the whole row (highlighted model) ---
$ws.Range("C6").Value = 0.925
$ws.Range("D6").Value = 0.102
$ws.Range("E6").Value = 5947.591
$ws.Range("A6:F6").Font.Bold = $true

# --- Row 10 (model 9): tweak CFI + AICc ---
$ws.Range("C10").Value = 0.896
$ws.Range("E10").Value = 3175.605

# --- Row 11 (model 10): now becomes "NA" across the board ---
$ws.Range("C11").Value = "NA"
$ws.Range("D11").Value = "NA"
$ws.Range("E11").Value = "NA"
$ws.Range("F11").Value = "NA"

# --- Row 12 (model 11): tweak CFI + AICc ---
$ws.Range("C12").Value = 0.932
$ws.Range("E12").Value = 5939.56

# --- Row 13 (model 12): now becomes "NA" across the board ---
$ws.Range("C13").Value = "NA"
$ws.Range("D13").Value = "NA"
$ws.Range("E13").Value = "NA"
$ws.Range("F13").Value = "NA"

# --- Selection moved to B12 on this sheet ---
[void]$ws.Range("B12").Select()
